# Updates to Rank and other scripts
#
# Re-scores several reaches' RTT_contaminants (column C) and the derived
# contaminants_dif (column D) on the EDT/RTT habitat attribute score
# comparison sheet for REACHES.
#
# The existing columns store their scores as text (e.g. "1", "4", "-2")
# rather than numbers, so each new value is entered with a leading "'"
# to force text entry (otherwise Excel would coerce the numeric-looking
# string into a Number cell). The style is reset back to Normal right
# after so the "number stored as text" quote-prefix formatting Excel
# applies to forced-text entries does not stick around on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextValue($rangeAddress, $value) {
    $rng = $ws.Range($rangeAddress)
    $rng.Value = "'" + $value
    $rng.Style = "Normal"
}

Set-TextValue "C5" "5"
Set-TextValue "D5" "0"
Set-TextValue "C6" "5"
Set-TextValue "D6" "0"
Set-TextValue "C11" "5"
Set-TextValue "D11" "0"
Set-TextValue "C16" "5"
Set-TextValue "D16" "-2"
Set-TextValue "C26" "5"
Set-TextValue "D26" "0"
Set-TextValue "C27" "5"
Set-TextValue "D27" "0"
Set-TextValue "C28" "5"
Set-TextValue "D28" "0"
Set-TextValue "C29" "5"
Set-TextValue "D29" "0"
Set-TextValue "C33" "5"
Set-TextValue "D33" "0"
Set-TextValue "C34" "5"
Set-TextValue "D34" "0"
Set-TextValue "C35" "5"
Set-TextValue "D35" "0"
Set-TextValue "C37" "5"
Set-TextValue "D37" "0"
Set-TextValue "C38" "5"
Set-TextValue "D38" "0"
Set-TextValue "C40" "5"
Set-TextValue "D40" "0"
Set-TextValue "C41" "5"
Set-TextValue "D41" "0"
Set-TextValue "C42" "5"
Set-TextValue "D42" "0"
Set-TextValue "C43" "5"
Set-TextValue "D43" "0"
Set-TextValue "C44" "5"
Set-TextValue "D44" "0"
Set-TextValue "C47" "5"
Set-TextValue "D47" "0"
Set-TextValue "C73" "5"
Set-TextValue "D73" "0"
Set-TextValue "C75" "5"
Set-TextValue "D75" "0"
Set-TextValue "C78" "5"
Set-TextValue "D78" "-2"
Set-TextValue "C81" "5"
Set-TextValue "D81" "0"
Set-TextValue "C101" "5"
Set-TextValue "D101" "0"
Set-TextValue "C111" "5"
Set-TextValue "D111" "0"
Set-TextValue "C112" "5"
Set-TextValue "D112" "0"
Set-TextValue "C115" "5"
Set-TextValue "D115" "0"
Set-TextValue "C118" "5"
Set-TextValue "D118" "0"
Set-TextValue "C119" "5"
Set-TextValue "D119" "0"
Set-TextValue "C125" "5"
Set-TextValue "C127" "5"
Set-TextValue "C128" "5"
